$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.092.65"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "3.592.19"
$ws.Range("E3").Value = "  -2.77%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "619.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -7.99%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "155.84"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.57%  "
$ws.Range("D7").Value = "3.587.35"
$ws.Range("E7").Value = "  -2.82%  "
$ws.Range("E8").Value = "  -0.01%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.486"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.87%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.141"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.62%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.03"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.29%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.431"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.94%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000223"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.82%  "
$ws.Range("D14").Value = "4.201.31"
$ws.Range("E14").Value = "  -2.74%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "31.69"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.92%  "
$ws.Range("D16").Value = "3.578.65"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "68.144.48"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("E18").Value = "  +0.42%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.38"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.33%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.51"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.18%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.83"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "453.86"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.14%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.637"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.26%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "77.65"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "3.736.13"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("E26").Value = "  +0.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.68"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.18%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0000116"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -9.61%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.36"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -8.66%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.57"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.66%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.62"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "25.97"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.91"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.22%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.592.81"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.159"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -5.90%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.40%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "8.11"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "177.12"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -6.53%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.60"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -8.77%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0874"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.25%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.899"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "45.88"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "28.50"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.80%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.97%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.20"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -7.84%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -6.76%  "
